$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new job posting row (Job_Id=13) at row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Nest JS Developer"
$ws.Range("C14").Value = "Hi"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
